$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "310.55"); Excel would
# auto-coerce these to real numbers on assignment, so force the
# range to Text format first, then restore the default "Normal"
# style afterwards so no stray number format sticks on the cells.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range('D2').Value = '43.485.36'
$ws.Range('E2').Value = '  +2.89%  '
$ws.Range('D3').Value = '2.312.39'
$ws.Range('E3').Value = '  +2.10%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '310.55'
$ws.Range('E5').Value = '  +0.65%  '
$ws.Range('D6').Value = '105.73'
$ws.Range('E6').Value = '  +8.74%  '
$ws.Range('E7').Value = '  +1.47%  '
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('E9').Value = '  +8.50%  '
$ws.Range('D10').Value = '37.01'
$ws.Range('E10').Value = '  +6.09%  '
$ws.Range('D11').Value = '52.90'
$ws.Range('E11').Value = '  +1.64%  '
$ws.Range('E12').Value = '  +1.19%  '
$ws.Range('E13').Value = '  -0.88%  '
$ws.Range('E14').Value = '  +3.11%  '
$ws.Range('D15').Value = '2.671.66'
$ws.Range('E15').Value = '  +2.05%  '
$ws.Range('D16').Value = '15.14'
$ws.Range('E16').Value = '  +4.08%  '
$ws.Range('D17').Value = '2.318.84'
$ws.Range('E17').Value = '  +1.76%  '
$ws.Range('E18').Value = '  +3.64%  '
$ws.Range('D19').Value = '43.391.62'
$ws.Range('E19').Value = '  +2.96%  '
$ws.Range('D20').Value = '12.22'
$ws.Range('E20').Value = '  -0.32%  '
$ws.Range('E21').Value = '  +3.22%  '
$ws.Range('D22').Value = '6.19'
$ws.Range('E22').Value = '  +4.24%  '
$ws.Range('D23').Value = '68.38'
$ws.Range('D24').Value = '242.64'
$ws.Range('E24').Value = '  +2.66%  '
$ws.Range('E25').Value = '  +3.71%  '
$ws.Range('E26').Value = '  +1.49%  '
$ws.Range('E27').Value = '  +0.09%  '
$ws.Range('D28').Value = '24.83'
$ws.Range('E28').Value = '  +5.56%  '
$ws.Range('E29').Value = '  +12.30%  '
$ws.Range('D30').Value = '37.32'
$ws.Range('E30').Value = '  +0.77%  '
$ws.Range('E31').Value = '  +1.37%  '
$ws.Range('D32').Value = '166.39'
$ws.Range('E32').Value = '  +1.88%  '
$ws.Range('E33').Value = '  +1.70%  '
$ws.Range('D34').Value = '18.45'
$ws.Range('E34').Value = '  +4.74%  '
$ws.Range('E36').Value = '  +6.73%  '
$ws.Range('D37').Value = '0.0745'
$ws.Range('E37').Value = '  +1.88%  '
$ws.Range('D38').Value = '3.06'
$ws.Range('E38').Value = '  -0.88%  '
$ws.Range('D39').Value = '4.55'
$ws.Range('E39').Value = '  +9.95%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').Value = '0.107'
$ws.Range('E40').Value = '  +3.32%  '
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').Value = '1.89'
$ws.Range('E41').Value = '  +4.28%  '
$ws.Range('B42').Value = 'ApeXProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D42').Value = '2.79'
$ws.Range('E42').Value = '  +22.56%  '
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').Value = '0.116'
$ws.Range('E43').Value = '  +1.12%  '
$ws.Range('E44').Value = '  +4.09%  '
$ws.Range('D45').Value = '1.996.87'
$ws.Range('E45').Value = '  +2.57%  '
$ws.Range('D46').Value = '3.15'
$ws.Range('E46').Value = '  +7.16%  '
$ws.Range('D47').Value = '19.00'
$ws.Range('E47').Value = '  +1.95%  '
$ws.Range('E48').Value = '  +2.75%  '
$ws.Range('D49').Value = '57.08'
$ws.Range('E49').Value = '  +5.13%  '
$ws.Range('E50').Value = '  +9.54%  '
$ws.Range('D51').Value = '2.92'
$ws.Range('E51').Value = '  +1.00%  '

$dRange.Style = "Normal"
